$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 39: "Save user feedback" service entry (Feedback form) ---
$ws.Range("B39").Value = "feedback"
$ws.Range("C39").Value = "Save user feedback"
$ws.Range("D39").Value = "WS-FED-02"
$ws.Range("E39").Value = "app.feedback.save"
$ws.Range("F39").Value = "'false"
$ws.Range("G39").Value = "feedback"
$ws.Range("H39").Value = "/save"
$ws.Range("I39").Value = "POST"

$ws.Range("M39").Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D39,"'',''CONNON_CONFIG'', 0, ''",C39,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Range("N39").Formula = '=_xlfn.CONCAT(IF(I39="GET","@GetMapping(",IF(I39="POST","@PostMapping(",IF(I39="DELETE","@DeleteMapping(",IF(I39="PUT","@PutMapping(","")))),CHAR(34),H39,CHAR(34),")")'
$ws.Range("O39").Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D39,,CHAR(34),", serviceName = ",CHAR(34),C39,CHAR(34), ", queryId = ",CHAR(34),E39,CHAR(34),", logActivity =",F39,")")'

# Match the row-38 cell formatting (borders/fills/number formats) for the new row
$ws.Range("B38:L38").Copy()
$ws.Range("B39:L39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Cursor moved on to the next empty row as part of this edit
$ws.Range("L43").Select() | Out-Null
